$d = $word.ActiveDocument

# Locate the specific TMPOS paragraph that needs the course correction
# (the "261202" position report, which incorrectly recorded course 057T
# instead of 000T).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "TMPOS/261202ZAPR*057T*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the TMPOS/261202 paragraph"
}

$r = $target.Range

# Find the "057T" token inside that paragraph and narrow down to the "057"
# digits that must become "000".
$found = $r.Find.Execute("057T", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate '057T' inside the target paragraph"
}

$courseStart = $r.Start
$courseEnd = $r.Start + 3

# Replace the course digits 057 -> 000.
$courseRange = $d.Range($courseStart, $courseEnd)
$courseRange.Text = "000"

# Word tracks the location of the last edit with the hidden "_GoBack"
# bookmark; move it to sit right after the corrected course value.
$goBackRange = $d.Range($courseEnd, $courseEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)
